$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 74 (quarter 01-01-2021) with revised figures
$ws.Range("B74").Value = -26743
$ws.Range("C74").Value = 441966
$ws.Range("D74").Value = 157871
$ws.Range("E74").Value = 199739
$ws.Range("F74").Value = 8790
$ws.Range("G74").Value = 35345
$ws.Range("H74").Value = 40220
$ws.Range("I74").Value = 468709
$ws.Range("J74").Value = 285787
$ws.Range("K74").Value = 116820
$ws.Range("L74").Value = 7853
$ws.Range("M74").Value = 58250

# Add new row 75 for quarter 01-04-2021
# A75 holds a date-like label ("01-04-2021") that must be stored as text
# (matching the other "Serie" column entries). Assigning the literal string
# directly gets auto-converted to a date serial by Excel, so build it as a
# text formula first and then convert the formula to its static text value.
$ws.Range("A75").Formula = "=""01-04-2021"""
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4163)

$ws.Range("B75").Value = -16701
$ws.Range("C75").Value = 454892
$ws.Range("D75").Value = 159638
$ws.Range("E75").Value = 201181
$ws.Range("F75").Value = 9083
$ws.Range("G75").Value = 40036
$ws.Range("H75").Value = 44954
$ws.Range("I75").Value = 471593
$ws.Range("J75").Value = 285942
$ws.Range("K75").Value = 120198
$ws.Range("L75").Value = 8333
$ws.Range("M75").Value = 57120
